$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

Set-TextValue "D2" "52.271.42"
Set-TextValue "E2" "  +6.02%  "
Set-TextValue "D3" "2.794.54"
Set-TextValue "E3" "  +6.54%  "
Set-TextValue "E4" "  -0.03%  "
Set-TextValue "D5" "116.64"
Set-TextValue "E5" "  +4.66%  "
Set-TextValue "D6" "340.76"
Set-TextValue "E6" "  +4.74%  "
Set-TextValue "D7" "0.552"
Set-TextValue "E7" "  +5.36%  "
Set-TextValue "D8" "0.999"
Set-TextValue "E8" "  -0.05%  "
Set-TextValue "E9" "  +6.21%  "
Set-TextValue "D10" "42.09"
Set-TextValue "E10" "  +7.01%  "
Set-TextValue "E11" "  +7.16%  "
Set-TextValue "D12" "20.11"
Set-TextValue "E12" "  +0.72%  "
Set-TextValue "E13" "  +2.60%  "
Set-TextValue "E14" "  +1.58%  "
Set-TextValue "D15" "3.236.47"
Set-TextValue "E15" "  +6.53%  "
Set-TextValue "D16" "2.802.69"
Set-TextValue "E16" "  +6.50%  "
Set-TextValue "E17" "  +4.23%  "
Set-TextValue "D18" "52.042.20"
Set-TextValue "E18" "  +5.55%  "
Set-TextValue "E19" "  +11.00%  "
Set-TextValue "D20" "13.35"
Set-TextValue "E20" "  +1.11%  "
Set-TextValue "E21" "  +4.68%  "
Set-TextValue "D22" "0.0₃0985"
Set-TextValue "E22" "  +4.21%  "
Set-TextValue "D23" "278.77"
Set-TextValue "E23" "  +4.17%  "
Set-TextValue "D24" "70.32"
Set-TextValue "E24" "  +2.04%  "
Set-TextValue "D25" "2.78"
Set-TextValue "E25" "  +8.68%  "
Set-TextValue "D26" "26.89"
Set-TextValue "E26" "  +3.76%  "
Set-TextValue "E27" "  +0.03%  "
Set-TextValue "D28" "10.26"
Set-TextValue "E28" "  +1.15%  "
Set-TextValue "E29" "  +1.28%  "
Set-TextValue "E30" "  +3.36%  "
Set-TextValue "D31" "34.85"
Set-TextValue "E31" "  +1.15%  "
Set-TextValue "D32" "50.42"
Set-TextValue "E32" "  +1.89%  "
Set-TextValue "E33" "  +5.48%  "
Set-TextValue "E34" "  +2.96%  "
Set-TextValue "E35" "  +5.27%  "
Set-TextValue "E36" "  -0.04%  "
Set-TextValue "D37" "18.98"
Set-TextValue "E37" "  +0.12%  "
Set-TextValue "D38" "4.99"
Set-TextValue "E38" "  +0.99%  "
Set-TextValue "D39" "3.24"
Set-TextValue "B40" "Stacks"
Set-TextValue "C40" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D40" "2.76"
Set-TextValue "E40" "  +29.78%  "
Set-TextValue "B41" "VeChain"
Set-TextValue "C41" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D41" "0.0377"
Set-TextValue "E41" "  +13.15%  "
Set-TextValue "B42" "Stellar"
Set-TextValue "C42" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D42" "0.116"
Set-TextValue "E42" "  +4.36%  "
Set-TextValue "B43" "EnergySwap"
Set-TextValue "C43" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D43" "23.37"
Set-TextValue "E43" "  +2.83%  "
Set-TextValue "E44" "  +2.96%  "
Set-TextValue "D45" "124.69"
Set-TextValue "E45" "  -3.34%  "
Set-TextValue "D46" "2.100.32"
Set-TextValue "E46" "  +2.06%  "
Set-TextValue "E47" "  +2.26%  "
Set-TextValue "E49" "  +7.62%  "
Set-TextValue "D50" "0.905"
Set-TextValue "E50" "  +22.22%  "
Set-TextValue "E51" "  +1.83%  "
